$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently sits inside the
#    "Troisième session : " paragraph, splitting it into two runs
#    ("Troisième" and " session : ") without the bookmark markers.
# ---------------------------------------------------------------------------
$troisiemePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Troisi*session*") {
        $troisiemePara = $p
        break
    }
}

$troisiemeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>Troisi&#232;me</w:t></w:r>' + `
              '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> session&#160;: </w:t></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$troisiemePara.Range.InsertXML($troisiemeXml)

# ---------------------------------------------------------------------------
# 2) Insert the new "Quatrième" .. "faire" paragraphs right after the
#    "Conjugaison du verbe être" paragraph (and before the trailing blank
#    paragraph). The very last paragraph ("faire") receives the "_GoBack"
#    bookmark that used to live in the "Troisième" paragraph.
# ---------------------------------------------------------------------------
$conjugaisonPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Conjugaison du verbe*") {
        $conjugaisonPara = $p
        break
    }
}

$insertRange = $conjugaisonPara.Next().Range
$insertRange.Collapse(1)

$newParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">Quatri&#232;me session&#160;: </w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>La conjugaison des verbes r&#233;guliers</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">Cinqui&#232;me session&#160;: </w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>Exercice des sessions pr&#233;c&#233;dentes</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>Huiti&#232;me session&#160;:</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>Les nombres 10-20</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>100</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>1000</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR" w:bidi="fa-IR"/></w:rPr><w:t>faire</w:t></w:r>' + `
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$insertRange.InsertXML($newParasXml)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
